$d = $word.ActiveDocument

# 1. Remove the "Admin can change state and submitter..." bullet entirely
#    (disallowing admins to change expense status is no longer a planned task).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Admin can change state and submitter*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# 2. Move the "Identify records which cannot be deleted or edited..." bullet
#    down to the end of the "High" priority section, right before the
#    "Disable edit button..." bullet.
$src = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Identify records which cannot be deleted or edited*") {
        $src = $p
        break
    }
}
if ($src -ne $null) {
    $src.Range.Cut()

    $dest = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "Disable edit button*") {
            $dest = $p
            break
        }
    }
    if ($dest -ne $null) {
        $insertionPoint = $d.Range($dest.Range.Start, $dest.Range.Start)
        $insertionPoint.Paste()
    }
}
